$wb = $excel.ActiveWorkbook

# --- Rename the first sheet ---
$wsModel = $wb.Worksheets.Item("Model_Karsilastirma")
$wsModel.Name = "Model_Sonuclari"

$wsHam = $wb.Worksheets.Item("Ham_Veri")

# --- Update header label: 3rd degree polynomial -> 2nd degree polynomial ---
$wsModel.Range("E1").Value = "Polinom (2.Der) (%)"

# --- Update the recalculated model values (Model_Sonuclari sheet) ---
$wsModel.Range("C2").Value = 11.99841028719179
$wsModel.Range("D2").Value = 28.51192875716237
$wsModel.Range("E2").Value = 28.29480949943364
$wsModel.Range("F2").Value = 28.98751674590322

$wsModel.Range("D3").Value = 26.45229895169499
$wsModel.Range("E3").Value = 28.93670644539078
$wsModel.Range("F3").Value = 27.61591402048322

$wsModel.Range("C4").Value = 25.37901045048939
$wsModel.Range("D4").Value = 24.32324342244781
$wsModel.Range("E4").Value = 27.24596210987605
$wsModel.Range("F4").Value = 24.89893254032222

$wsModel.Range("C5").Value = 22.33335939214919
$wsModel.Range("D5").Value = 22.91158703892523
$wsModel.Range("E5").Value = 24.80541519763349
$wsModel.Range("F5").Value = 23.05668112579334

$wsModel.Range("C6").Value = 19.28770833380899
$wsModel.Range("D6").Value = 21.49993065540264
$wsModel.Range("E6").Value = 21.31271005875625
$wsModel.Range("F6").Value = 21.36207068986936

$wsModel.Range("C7").Value = 16.44177209896651
$wsModel.Range("D7").Value = 20.18084190358646
$wsModel.Range("E7").Value = 17.09810541375526
$wsModel.Range("F7").Value = 19.95433944690444

# --- Update the raw data (Ham_Veri sheet) that feeds the model ---
$wsHam.Range("D2").Value = 33.57499999999958
$wsHam.Range("D6").Value = 25.52879656811268
$wsHam.Range("D7").Value = 16.44177209896651
